# Commit: "Changed creds for usser accounts"
# Rotates the sanitized test credentials (emails) used across the
# Network / Vendor 1 / Vendor 2 / Coseller sheets, and swaps the
# product URLs listed on the "List of Products" sheet (keeping only
# the first three, clearing the rest).

$wb = $excel.ActiveWorkbook

# ---- Network sheet -------------------------------------------------
$wsNetwork = $wb.Worksheets.Item("Network")
$wsNetwork.Range("A2").Value = "sanitynetwork@mailinator.com"

# ---- Vendor 1 sheet --------------------------------------------------
$wsVendor1 = $wb.Worksheets.Item("Vendor 1")
$wsVendor1.Range("A2").Value = "sanityvendor@mailinator.com"

# ---- Vendor 2 sheet --------------------------------------------------
$wsVendor2 = $wb.Worksheets.Item("Vendor 2")
$wsVendor2.Range("A2").Value = "sanityvendor1@mailinator.com"

# ---- Coseller sheet ---------------------------------------------------
$wsCoseller = $wb.Worksheets.Item("Coseller")
$wsCoseller.Range("A2").Value = "sanitycoseller@mailinator.com"

# ---- List of Products sheet -------------------------------------------
$wsProducts = $wb.Worksheets.Item("List of Products")
$wsProducts.Range("A2").Value = "https://beta.shoptype.com/product/5af42689-ef65-295b-0193-b31cdbd47c4b?tid=pv_71832512-5344-5d1a-3dae-5e5ebfec4658&utm_medium=copy"
$wsProducts.Range("A3").Value = "https://beta.shoptype.com/product/e9f70e52-589f-d095-b979-486d6fc5366d?tid=pv_295ecb00-15bc-0369-fc29-2731103beb41&utm_medium=copy"
$wsProducts.Range("A4").Value = "https://beta.shoptype.com/product/df98c477-7bc9-be70-2750-17dae5f0d360?tid=pv_7a6e2d1d-f6cd-4d70-73ec-e27ec59047de&utm_medium=copy"

# Row heights picked up by Excel when the longer URLs wrap onto the
# 75-wide "A" column.
$wsProducts.Rows.Item(2).RowHeight = 37.2
$wsProducts.Rows.Item(3).RowHeight = 30
$wsProducts.Rows.Item(4).RowHeight = 27

# Remaining product rows are no longer used.
$wsProducts.Range("A5").Value = $null
$wsProducts.Range("A6").Value = $null
$wsProducts.Range("A7").Value = $null

# ---- Selections / active sheet ----------------------------------------
# "Network" becomes the tab shown when the workbook is reopened
# (was "Vendor 1").
$null = $wsVendor1.Range("A3").Select()
$null = $wsVendor2.Range("A3").Select()
$null = $wsCoseller.Range("A3").Select()
$null = $wsProducts.Range("A5").Select()

$null = $wsNetwork.Activate()
$null = $wsNetwork.Range("G15").Select()
